# checking similarity for each action(keywords) instead of whole keywords or sentence
#
# The Attributes column (D) previously held the union of attributes for the
# whole sentence's matched use cases (one list literal concatenated per use
# case). It is switched to sets (Python repr of `set`) computed per
# individual action/keyword instead, which also changes which attributes
# show up. The Functionality column (C) drops extra use cases that are no
# longer considered matches once similarity is checked per action.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Functionality (C) updates: remove the second, no-longer-matching use case ---
$ws.Range("C5").Value  = "['Pay Bills ']"
$ws.Range("C7").Value  = "['Cheque Services']"
$ws.Range("C9").Value  = "['Limit Cash']"

# --- Attributes (D) updates: per-action attribute sets replace the old per-sentence lists ---
$ws.Range("D2").Value  = "{'Cust_Addr', 'Loan_Amt', 'Acc_num', 'Cus_Nme'}{'To_AcctNum', 'Debit_pin', 'Amt_avail', 'Acc_num', 'From_AcctNum'}"
$ws.Range("D3").Value  = "{'Debit_pin', 'Max_limit', 'Acc_num', 'Cus_Nme', 'Bill_type'}{'Cust_Addr', 'Loan_Amt', 'Acc_num', 'Cus_Nme'}"
$ws.Range("D4").Value  = "{'Amt_trnsfr', 'To_AcctNum', 'Amt_avail', 'Cus_Nme', 'From_AcctNum'}{'Cust_Addr', 'Loan_Amt', 'Acc_num', 'Cus_Nme'}"
$ws.Range("D5").Value  = "{'Debit_pin', 'Amt_avail', 'Acc_num', 'Bill_type', 'From_AcctNum'}"
$ws.Range("D6").Value  = "{'Debit_pin', 'To_AcctNum', 'Acc_num', 'Cus_Nme', 'From_AcctNum'}{'Cred_Score', 'Loan_Amt', 'Loan_purp'}"
$ws.Range("D7").Value  = "{'Cust_Addr', 'Loan_Amt', 'Acc_num', 'Cus_Nme'}"
$ws.Range("D8").Value  = "{'Acc_type', 'Max_limit', 'Acc_num', 'Cus_Nme'}{'Acc_type', 'Debit_pin', 'Amt_deposit', 'Acc_num', 'Cus_Nme'}"
$ws.Range("D9").Value  = "{'Acc_type', 'Debit_pin', 'Amt_deposit', 'Acc_num', 'Cus_Nme'}"
$ws.Range("D10").Value = "{'Acc_num'}{'Cust_Addr', 'Loan_Amt', 'Acc_num', 'Cus_Nme'}"
$ws.Range("D11").Value = "{'Cred_Score', 'Loan_Amt', 'Loan_purp'}{'Acc_num'},{'Debit_pin', 'To_AcctNum', 'Acc_num', 'Cus_Nme', 'From_AcctNum'}{'Cred_Score', 'Loan_Amt', 'Loan_purp'}"

# D12 had no real content (empty string cell) and is removed entirely.
$ws.Range("D12").ClearContents()
